$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-27 Saturday" "2025-09-28 Sunday"

Replace-Text "271÷9=" "983÷5="
Replace-Text "198÷2=" "399÷2="
Replace-Text "164÷3=" "560÷9="
Replace-Text "243÷2=" "917÷8="
Replace-Text "731÷3=" "644÷8="
Replace-Text "635÷8=" "985÷5="
Replace-Text "355÷9=" "500÷4="
Replace-Text "565÷7=" "874÷2="
Replace-Text "325÷9=" "934÷2="
Replace-Text "148÷3=" "322÷5="
Replace-Text "593÷6=" "975÷3="
Replace-Text "935÷5=" "504÷8="
Replace-Text "579÷8=" "319÷5="
Replace-Text "490÷9=" "776÷7="
Replace-Text "172÷4=" "627÷2="
Replace-Text "366÷6=" "179÷6="
Replace-Text "815÷2=" "757÷9="
Replace-Text "188÷7=" "876÷5="
Replace-Text "263÷5=" "930÷3="
Replace-Text "148÷6=" "417÷2="
Replace-Text "217÷3=" "798÷5="
Replace-Text "352÷8=" "860÷3="
Replace-Text "108÷2=" "733÷8="
Replace-Text "695÷8=" "280÷3="
Replace-Text "297÷8=" "965÷9="
